# Adds a new "Find Peak Element" (LeetCode 162) entry as row 7 of the
# "数组" (Array) worksheet.
#
# Columns: A=No. | B=leetcode # | C=题目(problem) | D=解题方法(solution)
#          | E=解题关键词(keywords) | F=时间复杂度(time) | G=空间复杂度(space)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("数组")

$solution = "1 基于二分查找思想，获取数组中间节点mid`n2 num[mid] < num[mid+1] 升序，峰值在右侧，mid肯定不是峰值，low = mid+1`n3 num[mid] > num[mid+1] 这部分是降序，峰值在左侧，有可能mid就是峰值，high = mid - 1`n4 不会出现相邻元素相等的情况，代码中可以处理这个情况`n5 最终可以返回low/high索引"

$problem = "峰值元素是指其值大于左右相邻值的元素。 `n 给定一个输入数组 nums，其中 nums[i] ≠ nums[i+1]，找到峰值元素并返回其索引。 `n 数组可能包含多个峰值，在这种情况下，返回任何一个峰值所在位置即可。 `n 你可以假设 nums[-1] = nums[n] = -∞。 `n 示例 1: `n 输入: nums = [1,2,3,1]`n输出: 2`n解释: 3 是峰值元素，你的函数应该返回其索引 2。 `n 示例 2: `n 输入: nums = [1,2,1,3,5,6,4]`n输出: 1 或 5 `n解释: 你的函数可以返回索引 1，其峰值元素为 2；`n     或者返回索引 5， 其峰值元素为 6。`n 说明: `n 你的解法应该是 O(logN) 时间复杂度的。 `n Related Topics 数组 二分查找"

$keywords = "二分查找`n升序/降序`n"

# Write the new string-valued cells in the same order the original edit
# created them in, so new shared-string entries line up (solution text,
# then the problem text, then the keywords).
$ws.Cells.Item(7, 4).Value = $solution
$ws.Cells.Item(7, 3).Value = $problem
$ws.Cells.Item(7, 5).Value = $keywords

$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 162
$ws.Cells.Item(7, 6).Value = "O(logN)"
$ws.Cells.Item(7, 7).Value = "O(1)"

# The row holding this (very long) entry is given the same "maxed out"
# row height used by the sheet's other long entries.
$ws.Rows.Item(7).RowHeight = 409.6

# Leave the new row's first description cell selected, as in the source
# edit.
$ws.Activate()
$ws.Range("D7").Select() | Out-Null
